$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: rewrite the "SELECT name_author, sum(amount) AS " paragraph so the
# SQL keywords/identifiers are split across many runs (some with the en-US
# language mark, some without) instead of a few proofed runs, and drop the
# paragraph-mark rPr (en-US lang) from pPr. The trailing "Количество" run is
# recreated identically.
# ---------------------------------------------------------------------------
$oldFrag1 = "SELECT name_author, sum(amount) AS "
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq ($oldFrag1 + "Количество`r")) {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:after="0"/><w:ind w:left="1800"/><w:jc w:val="both"/></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>SELECT</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>name</w:t></w:r><w:r><w:t>_</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>author</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>sum</w:t></w:r><w:r><w:t>(</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>amount</w:t></w:r><w:r><w:t xml:space="preserve">) </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>AS</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Количество</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    [void]$target.Range.InsertXML($xml1)
}

# ---------------------------------------------------------------------------
# Edit 2: append a page break plus a brand-new "Задание" task (with its SQL
# answer) at the end of the document, right before the final section break.
# ---------------------------------------------------------------------------
$xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:spacing w:line="259" w:lineRule="auto"/></w:pPr><w:r><w:br w:type="page"/></w:r></w:p><w:p><w:pPr><w:spacing w:after="0"/><w:ind w:left="1800"/><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:lastRenderedPageBreak/><w:t>Задание</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0"/><w:ind w:left="1800"/><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">Если в таблицах </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>supply</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">  и</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>book</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> есть одинаковые книги, которые имеют равную цену,  вывести их название и автора, а также посчитать общее количество экземпляров книг в таблицах </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>supply</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> и </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>book</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>,  столбцы назвать Название, Автор  и Количество.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0"/><w:ind w:left="1800"/><w:jc w:val="both"/></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0"/><w:ind w:left="1800"/><w:jc w:val="both"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">SELECT </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>book.title</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> AS </w:t></w:r><w:r><w:t>Название</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>author.name_author</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> AS </w:t></w:r><w:r><w:t>Автор</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>supply.amount</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> + </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>book.amount</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> AS </w:t></w:r><w:r><w:t>Количество</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0"/><w:ind w:left="1800"/><w:jc w:val="both"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>FROM author</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0"/><w:ind w:left="1800"/><w:jc w:val="both"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">    INNER JOIN book </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>USING(</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>author_id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>)</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0"/><w:ind w:left="1800"/><w:jc w:val="both"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">    INNER JOIN supply</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0"/><w:ind w:left="1800"/><w:jc w:val="both"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">    ON </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>book.title</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> = </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>supply.title</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>book.price</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> = </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>supply.price</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>;</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:after="0"/><w:ind w:left="1800"/><w:jc w:val="both"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0"/><w:ind w:left="1800"/><w:jc w:val="both"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:spacing w:after="0"/><w:ind w:left="1800"/><w:jc w:val="both"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)
[void]$endRange.InsertXML($xml2)

Write-Output "edit complete"
